# Apply the ModeloRelacional.docx edit: rewrite the trailing functional-dependency
# list of the "Enfermeiro" paragraph from
#   DiaTrabalho-> DiaTrabalho->DiaAbsoluto ->Departamento
# to
#   DiaAbsoluto->DiaTrabalho, Designação->Departamento
# and relocate the _GoBack bookmark so it still sits right after "DiaAbsoluto"
# (now split as "Dia" + "Absoluto"), matching the target run layout:
#   "Dia" | "Absoluto" | [_GoBack] | "->" | "DiaTrabalho" | ", Designação" | "->" | "Departamento" | ")"

$d = $word.ActiveDocument

# --- locate the paragraph that still has the old (unedited) wording ---------
$oldFragment = "DiaTrabalho-> DiaTrabalho->DiaAbsoluto ->Departamento"
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($oldFragment)) {
        $para = $p
    }
}

$paraStart  = $para.Range.Start
$paraText   = $para.Range.Text
$fragStart  = $paraStart + $paraText.IndexOf($oldFragment)
$fragEnd    = $fragStart + $oldFragment.Length

# --- rewrite the fragment's text in one shot --------------------------------
$newFragment = "DiaAbsoluto->DiaTrabalho, Designação->Departamento"
$target = $d.Range($fragStart, $fragEnd)
$target.Text = $newFragment

# --- work out the absolute offsets of each logical piece in the new text ---
$pieces = @("Dia", "Absoluto", "->", "DiaTrabalho", ", Designação", "->", "Departamento")

$boundaries = @()
$boundaries += $fragStart               # split before "Dia" (keeps preceding ", " run intact)
$pos = $fragStart
foreach ($piece in $pieces) {
    $pos = $pos + $piece.Length
    $boundaries += $pos                 # split after each piece (last one precedes ")")
}

# --- temporarily bookmark every boundary so the engine can't re-merge runs -
$tempNames = @()
$i = 0
foreach ($b in $boundaries) {
    $name = "TmpSplit" + $i
    $d.Bookmarks.Add($name, $d.Range($b, $b))
    $tempNames += $name
    $i = $i + 1
}

# --- move _GoBack to sit right after "Absoluto" (i.e. after "DiaAbsoluto") -
$gobackPos = $fragStart + "Dia".Length + "Absoluto".Length
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($gobackPos, $gobackPos))

# --- drop the scaffolding bookmarks, the run split they enforced persists --
foreach ($name in $tempNames) {
    $d.Bookmarks.Item($name).Delete()
}

Write-Output $para.Range.Text
